# The published page footer ("Ver no Jupiter ... / © 2020 ...") was dropped
# from the rebuilt site, along with the blank line that separated it from
# the bibliography. Remove that trailing block, keeping everything from the
# last bibliography entry ("FLEMMING, ...") onward unchanged, and leaving the
# blank paragraph + page-break paragraph that originally followed the footer.

$d = $word.ActiveDocument

# Locate the last bibliography paragraph that must be kept.
$flemRange = $d.Content
$flemRange.Find.Execute("FLEMMING, Diva M.; GONÇALVES, Mirian B. Cálculo A.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$flemIndex = $flemRange.Paragraphs.First.Index

# Locate the final line of the footer block that must be removed.
$copyRange = $d.Content
$copyRange.Find.Execute("© 2020 . Contact: luizeleno@usp.br", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$copyIndex = $copyRange.Paragraphs.First.Index

# Delete the blank paragraph + "Ver no Jupiter ..." + "© 2020 ..." paragraphs
# as a single range, which merges their paragraph marks away cleanly.
$startOfBlock = $d.Paragraphs.Item($flemIndex + 1).Range.Start
$endOfBlock = $d.Paragraphs.Item($copyIndex).Range.End

$d.Range($startOfBlock, $endOfBlock).Delete()
